$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"100.8744713333333"
$ws.Range("H2").Value = [double]"302.623414"
$ws.Range("I2").Value = [double]"0.1452075237922473"
$ws.Range("J2").Value = [double]"0.1452075237922473"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"0.1863256666666667"
$ws.Range("N2").Value = [double]"0.5589770000000001"
$ws.Range("O2").Value = [double]"0.01657678358851065"
$ws.Range("P2").Value = [double]"0.01657678358851065"
$ws.Range("Q2").Value = [double]"18.79550312083089"
$ws.Range("R2").Value = [double]"169.159528087478"
$ws.Range("S2").Value = [double]"0.002407073697327595"
$ws.Range("T2").Value = [double]"0.002407073697327595"
$ws.Range("G3").Value = [double]"100.8744713333333"
$ws.Range("H3").Value = [double]"302.623414"
$ws.Range("I3").Value = [double]"0.1452075237922473"
$ws.Range("J3").Value = [double]"0.1452075237922473"
$ws.Range("O3").Value = [double]"0.5186672939413604"
$ws.Range("P3").Value = [double]"0.5186672939413604"
$ws.Range("Q3").Value = [double]"588.0883157999665"
$ws.Range("R3").Value = [double]"5292.794842199698"
$ws.Range("S3").Value = [double]"0.07531439342525062"
$ws.Range("T3").Value = [double]"0.07531439342525062"
$ws.Range("G4").Value = [double]"100.8744713333333"
$ws.Range("H4").Value = [double]"302.623414"
$ws.Range("I4").Value = [double]"0.1452075237922473"
$ws.Range("J4").Value = [double]"0.1452075237922473"
$ws.Range("M4").Value = [double]"5.212463666666667"
$ws.Range("N4").Value = [double]"15.637391"
$ws.Range("O4").Value = [double]"0.4637358003923671"
$ws.Range("P4").Value = [double]"0.4637358003923669"
$ws.Range("Q4").Value = [double]"525.8045167192083"
$ws.Range("R4").Value = [double]"4732.240650472874"
$ws.Range("S4").Value = [double]"0.06733792726879149"
$ws.Range("T4").Value = [double]"0.06733792726879148"
$ws.Range("G5").Value = [double]"100.8744713333333"
$ws.Range("H5").Value = [double]"302.623414"
$ws.Range("I5").Value = [double]"0.1452075237922473"
$ws.Range("J5").Value = [double]"0.1452075237922473"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.01146633333333333"
$ws.Range("N5").Value = [double]"0.034399"
$ws.Range("O5").Value = [double]"0.001020122077762015"
$ws.Range("P5").Value = [double]"0.001020122077762015"
$ws.Range("Q5").Value = [double]"1.156660313131778"
$ws.Range("R5").Value = [double]"10.409942818186"
$ws.Range("S5").Value = [double]"0.0001481294008776246"
$ws.Range("T5").Value = [double]"0.0001481294008776245"
$ws.Range("I6").Value = [double]"0.7769829249672668"
$ws.Range("J6").Value = [double]"0.776982924967267"
$ws.Range("K6").Value = [double]"3"
$ws.Range("L6").Value = [double]"1"
$ws.Range("M6").Value = [double]"0.1863256666666667"
$ws.Range("N6").Value = [double]"0.5589770000000001"
$ws.Range("O6").Value = [double]"0.01657678358851065"
$ws.Range("P6").Value = [double]"0.01657678358851065"
$ws.Range("Q6").Value = [double]"100.5718203138609"
$ws.Range("R6").Value = [double]"905.146382824748"
$ws.Range("S6").Value = [double]"0.01287987779915039"
$ws.Range("T6").Value = [double]"0.01287987779915039"
$ws.Range("I7").Value = [double]"0.7769829249672668"
$ws.Range("J7").Value = [double]"0.776982924967267"
$ws.Range("O7").Value = [double]"0.5186672939413604"
$ws.Range("P7").Value = [double]"0.5186672939413604"
$ws.Range("S7").Value = [double]"0.4029956311314153"
$ws.Range("T7").Value = [double]"0.4029956311314154"
$ws.Range("I8").Value = [double]"0.7769829249672668"
$ws.Range("J8").Value = [double]"0.776982924967267"
$ws.Range("M8").Value = [double]"5.212463666666667"
$ws.Range("N8").Value = [double]"15.637391"
$ws.Range("O8").Value = [double]"0.4637358003923671"
$ws.Range("P8").Value = [double]"0.4637358003923669"
$ws.Range("Q8").Value = [double]"2813.498369037698"
$ws.Range("R8").Value = [double]"25321.48532133929"
$ws.Range("S8").Value = [double]"0.360314798600898"
$ws.Range("T8").Value = [double]"0.3603147986008979"
$ws.Range("I9").Value = [double]"0.7769829249672668"
$ws.Range("J9").Value = [double]"0.776982924967267"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.01146633333333333"
$ws.Range("N9").Value = [double]"0.034399"
$ws.Range("O9").Value = [double]"0.001020122077762015"
$ws.Range("P9").Value = [double]"0.001020122077762015"
$ws.Range("Q9").Value = [double]"6.189109832741778"
$ws.Range("R9").Value = [double]"55.70198849467599"
$ws.Range("S9").Value = [double]"0.000792617435803216"
$ws.Range("T9").Value = [double]"0.000792617435803216"
$ws.Range("G10").Value = [double]"53.798087"
$ws.Range("H10").Value = [double]"161.394261"
$ws.Range("I10").Value = [double]"0.07744166482137986"
$ws.Range("J10").Value = [double]"0.07744166482137986"
$ws.Range("K10").Value = [double]"3"
$ws.Range("L10").Value = [double]"1"
$ws.Range("M10").Value = [double]"0.1863256666666667"
$ws.Range("N10").Value = [double]"0.5589770000000001"
$ws.Range("O10").Value = [double]"0.01657678358851065"
$ws.Range("P10").Value = [double]"0.01657678358851065"
$ws.Range("Q10").Value = [double]"10.02396442566634"
$ws.Range("R10").Value = [double]"90.215679830997"
$ws.Range("S10").Value = [double]"0.001283733718477992"
$ws.Range("T10").Value = [double]"0.001283733718477992"
$ws.Range("G11").Value = [double]"53.798087"
$ws.Range("H11").Value = [double]"161.394261"
$ws.Range("I11").Value = [double]"0.07744166482137986"
$ws.Range("J11").Value = [double]"0.07744166482137986"
$ws.Range("O11").Value = [double]"0.5186672939413604"
$ws.Range("P11").Value = [double]"0.5186672939413604"
$ws.Range("Q11").Value = [double]"313.6375929301697"
$ws.Range("R11").Value = [double]"2822.738336371527"
$ws.Range("S11").Value = [double]"0.04016645873121893"
$ws.Range("T11").Value = [double]"0.04016645873121893"
$ws.Range("G12").Value = [double]"53.798087"
$ws.Range("H12").Value = [double]"161.394261"
$ws.Range("I12").Value = [double]"0.07744166482137986"
$ws.Range("J12").Value = [double]"0.07744166482137986"
$ws.Range("M12").Value = [double]"5.212463666666667"
$ws.Range("N12").Value = [double]"15.637391"
$ws.Range("O12").Value = [double]"0.4637358003923671"
$ws.Range("P12").Value = [double]"0.4637358003923669"
$ws.Range("Q12").Value = [double]"280.4205738236724"
$ws.Range("R12").Value = [double]"2523.785164413051"
$ws.Range("S12").Value = [double]"0.03591247241966"
$ws.Range("T12").Value = [double]"0.03591247241966"
$ws.Range("G13").Value = [double]"53.798087"
$ws.Range("H13").Value = [double]"161.394261"
$ws.Range("I13").Value = [double]"0.07744166482137986"
$ws.Range("J13").Value = [double]"0.07744166482137986"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.01146633333333333"
$ws.Range("N13").Value = [double]"0.034399"
$ws.Range("O13").Value = [double]"0.001020122077762015"
$ws.Range("P13").Value = [double]"0.001020122077762015"
$ws.Range("Q13").Value = [double]"0.6168667982376667"
$ws.Range("R13").Value = [double]"5.551801184139"
$ws.Range("S13").Value = [double]"7.899995202293556E-05"
$ws.Range("T13").Value = [double]"7.899995202293555E-05"
$ws.Range("E14").Value = [double]"3"
$ws.Range("F14").Value = [double]"1"
$ws.Range("G14").Value = [double]"0.2555676666666666"
$ws.Range("H14").Value = [double]"0.7667029999999999"
$ws.Range("I14").Value = [double]"0.0003678864191059829"
$ws.Range("J14").Value = [double]"0.000367886419105983"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"0.1863256666666667"
$ws.Range("N14").Value = [double]"0.5589770000000001"
$ws.Range("O14").Value = [double]"0.01657678358851065"
$ws.Range("P14").Value = [double]"0.01657678358851065"
$ws.Range("Q14").Value = [double]"0.04761881587011111"
$ws.Range("R14").Value = [double]"0.428569342831"
$ws.Range("S14").Value = [double]"6.098373554672008E-06"
$ws.Range("T14").Value = [double]"6.098373554672008E-06"
$ws.Range("E15").Value = [double]"3"
$ws.Range("F15").Value = [double]"1"
$ws.Range("G15").Value = [double]"0.2555676666666666"
$ws.Range("H15").Value = [double]"0.7667029999999999"
$ws.Range("I15").Value = [double]"0.0003678864191059829"
$ws.Range("J15").Value = [double]"0.000367886419105983"
$ws.Range("O15").Value = [double]"0.5186672939413604"
$ws.Range("P15").Value = [double]"0.5186672939413604"
$ws.Range("Q15").Value = [double]"1.489934536224555"
$ws.Range("R15").Value = [double]"13.409410826021"
$ws.Range("S15").Value = [double]"0.0001908106534754773"
$ws.Range("T15").Value = [double]"0.0001908106534754774"
$ws.Range("E16").Value = [double]"3"
$ws.Range("F16").Value = [double]"1"
$ws.Range("G16").Value = [double]"0.2555676666666666"
$ws.Range("H16").Value = [double]"0.7667029999999999"
$ws.Range("I16").Value = [double]"0.0003678864191059829"
$ws.Range("J16").Value = [double]"0.000367886419105983"
$ws.Range("M16").Value = [double]"5.212463666666667"
$ws.Range("N16").Value = [double]"15.637391"
$ws.Range("O16").Value = [double]"0.4637358003923671"
$ws.Range("P16").Value = [double]"0.4637358003923669"
$ws.Range("Q16").Value = [double]"1.332137176874778"
$ws.Range("R16").Value = [double]"11.989234591873"
$ws.Range("S16").Value = [double]"0.0001706021030175948"
$ws.Range("T16").Value = [double]"0.0001706021030175948"
$ws.Range("E17").Value = [double]"3"
$ws.Range("F17").Value = [double]"1"
$ws.Range("G17").Value = [double]"0.2555676666666666"
$ws.Range("H17").Value = [double]"0.7667029999999999"
$ws.Range("I17").Value = [double]"0.0003678864191059829"
$ws.Range("J17").Value = [double]"0.000367886419105983"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.01146633333333333"
$ws.Range("N17").Value = [double]"0.034399"
$ws.Range("O17").Value = [double]"0.001020122077762015"
$ws.Range("P17").Value = [double]"0.001020122077762015"
$ws.Range("Q17").Value = [double]"0.002930424055222222"
$ws.Range("R17").Value = [double]"0.02637381649699999"
$ws.Range("S17").Value = [double]"3.752890582388227E-07"
$ws.Range("T17").Value = [double]"3.752890582388227E-07"
